$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.302.53"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.929.75"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'248.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'0.7165"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.3209"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.22%  "
$ws.Range("D9").Value = "'27.68"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.71%  "
$ws.Range("D10").Value = "'0.07103"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").Value = "'0.7913"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.08000"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "1.931.23"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'5.394"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").Value = "'94.90"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "30.303.68"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'257.48"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "'0.000008068"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.84%  "
$ws.Range("D20").Value = "'5.773"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").Value = "2.182.70"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'6.835"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("D26").Value = "'164.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("E28").Value = "  -6.19%  "
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("D30").Value = "'1.354"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("D32").Value = "'4.399"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("D33").Value = "'4.143"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "'1.267"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").Value = "'0.7445"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").Value = "'2.767"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").Value = "'0.01966"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").Value = "'78.63"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("D41").Value = "'6.363"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.80%  "
$ws.Range("D42").Value = "'0.4516"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").Value = "'1.999"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'100.58"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("D47").Value = "'9.786"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "'7.441"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'36.74"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06107"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.4204"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.87%  "
